$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = '''71.164.79'
$ws.Range("E2").Value = '  +2.49%  '

# Row 3
$ws.Range("D3").Value = '''3.645.19'
$ws.Range("E3").Value = '  +4.22%  '

# Row 4
$ws.Range("D4").Value = '''0.998'
$ws.Range("E4").Value = '  +0.05%  '

# Row 5
$ws.Range("D5").Value = '''606.30'
$ws.Range("E5").Value = '  +1.26%  '

# Row 6
$ws.Range("D6").Value = '''202.66'
$ws.Range("E6").Value = '  +4.55%  '

# Row 7
$ws.Range("E7").Value = '  +1.61%  '

# Row 8
$ws.Range("D8").Value = '''1.00'
$ws.Range("E8").Value = '  +0.08%  '

# Row 9
$ws.Range("E9").Value = '  +10.76%  '

# Row 10
$ws.Range("D10").Value = '''0.650'
$ws.Range("E10").Value = '  +1.22%  '

# Row 11
$ws.Range("D11").Value = '''53.94'
$ws.Range("E11").Value = '  +2.01%  '

# Row 12
$ws.Range("E12").Value = '  +3.37%  '

# Row 13
$ws.Range("D13").Value = '''9.63'
$ws.Range("E13").Value = '  +2.39%  '

# Row 14
$ws.Range("D14").Value = '''4.221.65'
$ws.Range("E14").Value = '  +4.03%  '

# Row 15
$ws.Range("D15").Value = '''678.25'
$ws.Range("E15").Value = '  +13.94%  '

# Row 16
$ws.Range("D16").Value = '''71.120.89'
$ws.Range("E16").Value = '  +2.22%  '

# Row 17
$ws.Range("E17").Value = '  +2.82%  '

# Row 18
$ws.Range("D18").Value = '''3.643.89'
$ws.Range("E18").Value = '  +4.19%  '

# Row 19
$ws.Range("D19").Value = '''19.11'
$ws.Range("E19").Value = '  +1.23%  '

# Row 20
$ws.Range("E20").Value = '  +0.39%  '

# Row 21
$ws.Range("E21").Value = '  +2.63%  '

# Row 22
$ws.Range("D22").Value = '''18.98'
$ws.Range("E22").Value = '  +6.98%  '

# Row 23
$ws.Range("D23").Value = '''5.41'
$ws.Range("E23").Value = '  +2.88%  '

# Row 24
$ws.Range("D24").Value = '''105.33'
$ws.Range("E24").Value = '  +3.10%  '

# Row 25
$ws.Range("E25").Value = '  +1.34%  '

# Row 26
$ws.Range("E26").Value = '  -1.41%  '

# Row 27
$ws.Range("D27").Value = '''10.63'
$ws.Range("E27").Value = '  -1.09%  '

# Row 28
$ws.Range("D28").Value = '''9.96'
$ws.Range("E28").Value = '  +5.60%  '

# Row 29
$ws.Range("D29").Value = '''34.46'
$ws.Range("E29").Value = '  +4.88%  '

# Row 30
$ws.Range("D30").Value = '''4.57'
$ws.Range("E30").Value = '  +7.25%  '

# Row 31
$ws.Range("D31").Value = '''7.29'
$ws.Range("E31").Value = '  +5.03%  '

# Row 32
$ws.Range("D32").Value = '''12.27'
$ws.Range("E32").Value = '  +0.29%  '

# Row 33
$ws.Range("E33").Value = '  +2.01%  '

# Row 34
$ws.Range("D34").Value = '''63.42'
$ws.Range("E34").Value = '  +0.57%  '

# Row 35
$ws.Range("D35").Value = '''0.0₃0878'
$ws.Range("E35").Value = '  +9.01%  '

# Row 36
$ws.Range("D36").Value = '''3.924.11'
$ws.Range("E36").Value = '  +5.00%  '

# Row 37
$ws.Range("B37").Value = 'Bittensor'
$ws.Range("C37").Value = 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao'
$ws.Range("D37").Value = '''521.46'
$ws.Range("E37").Value = '  +5.77%  '

# Row 38
$ws.Range("B38").Value = 'Dai'
$ws.Range("C38").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D38").Value = '''1.00'
$ws.Range("E38").Value = '  +0.02%  '

# Row 39
$ws.Range("D39").Value = '''3.05'
$ws.Range("E39").Value = '  -3.95%  '

# Row 40
$ws.Range("D40").Value = '''3.63'
$ws.Range("E40").Value = '  -0.44%  '

# Row 41
$ws.Range("D41").Value = '''0.392'
$ws.Range("E41").Value = '  +1.46%  '

# Row 42
$ws.Range("D42").Value = '''36.74'
$ws.Range("E42").Value = '  +2.50%  '

# Row 43
$ws.Range("E43").Value = '  +3.69%  '

# Row 44
$ws.Range("E44").Value = '  +2.93%  '

# Row 45
$ws.Range("E45").Value = '  +9.86%  '

# Row 46
$ws.Range("D46").Value = '''3.46'
$ws.Range("E46").Value = '  +7.78%  '

# Row 47
$ws.Range("E47").Value = '  +1.63%  '

# Row 48
$ws.Range("D48").Value = '''8.68'
$ws.Range("E48").Value = '  +3.76%  '

# Row 49
$ws.Range("E49").Value = '  -0.29%  '

# Row 50
$ws.Range("D50").Value = '''0.000249'
$ws.Range("E50").Value = '  +2.51%  '

# Row 51
$ws.Range("E51").Value = '  +4.51%  '

